# Refresh the cryptos list (price + 1h volume columns) with the latest
# scrape, as produced by the "Updated cryptos list ... with GitHub Actions"
# automation. Also reflects Hedera/MXToken swapping rank positions 38/39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation (avoid Excel auto-converting numeric-looking
# strings like "0.8900" or "1.080" into numbers and losing the exact
# formatting), then reset the style back to Normal so no stray
# NumberFormat/style index is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.378.87"
Set-TextValue $ws.Range("E2") "  +1.27%  "

Set-TextValue $ws.Range("D3") "1.861.46"
Set-TextValue $ws.Range("E3") "  +1.84%  "

Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  -0.43%  "

Set-TextValue $ws.Range("D5") "315.37"
Set-TextValue $ws.Range("E5") "  +0.98%  "

Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  -0.42%  "

Set-TextValue $ws.Range("D7") "0.4618"
Set-TextValue $ws.Range("E7") "  +0.25%  "

Set-TextValue $ws.Range("D8") "0.3716"
Set-TextValue $ws.Range("E8") "  +0.55%  "

Set-TextValue $ws.Range("D9") "0.07318"
Set-TextValue $ws.Range("E9") "  -0.39%  "

Set-TextValue $ws.Range("D10") "0.8890"
Set-TextValue $ws.Range("E10") "  +2.08%  "

Set-TextValue $ws.Range("D11") "20.03"
Set-TextValue $ws.Range("E11") "  +1.25%  "

Set-TextValue $ws.Range("D12") "0.07834"
Set-TextValue $ws.Range("E12") "  -1.25%  "

Set-TextValue $ws.Range("D13") "1.919.68"
Set-TextValue $ws.Range("E13") "  +5.31%  "

Set-TextValue $ws.Range("D14") "5.395"
Set-TextValue $ws.Range("E14") "  +0.99%  "

Set-TextValue $ws.Range("D15") "6.552"
Set-TextValue $ws.Range("E15") "  +0.27%  "

Set-TextValue $ws.Range("D16") "91.78"
Set-TextValue $ws.Range("E16") "  +0.12%  "

Set-TextValue $ws.Range("E17") "  -0.48%  "

Set-TextValue $ws.Range("D18") "0.000008959"
Set-TextValue $ws.Range("E18") "  +0.95%  "

Set-TextValue $ws.Range("D19") "1.002"
Set-TextValue $ws.Range("E19") "  -0.41%  "

Set-TextValue $ws.Range("D20") "14.80"
Set-TextValue $ws.Range("E20") "  +0.78%  "

Set-TextValue $ws.Range("D21") "27.389.78"
Set-TextValue $ws.Range("E21") "  +2.41%  "

Set-TextValue $ws.Range("D22") "5.132"
Set-TextValue $ws.Range("E22") "  +0.30%  "

Set-TextValue $ws.Range("D23") "10.56"
Set-TextValue $ws.Range("E23") "  -0.02%  "

Set-TextValue $ws.Range("D24") "2.051.66"
Set-TextValue $ws.Range("E24") "  +5.29%  "

Set-TextValue $ws.Range("D25") "1.939"
Set-TextValue $ws.Range("E25") "  +4.90%  "

Set-TextValue $ws.Range("D26") "152.05"
Set-TextValue $ws.Range("E26") "  -0.21%  "

Set-TextValue $ws.Range("D27") "18.41"
Set-TextValue $ws.Range("E27") "  -0.43%  "

Set-TextValue $ws.Range("D28") "2.053"
Set-TextValue $ws.Range("E28") "  -0.76%  "

Set-TextValue $ws.Range("D29") "5.102"
Set-TextValue $ws.Range("E29") "  +0.28%  "

Set-TextValue $ws.Range("D30") "116.27"
Set-TextValue $ws.Range("E30") "  +0.86%  "

Set-TextValue $ws.Range("D31") "0.08844"
Set-TextValue $ws.Range("E31") "  -0.24%  "

Set-TextValue $ws.Range("D32") "3.110"
Set-TextValue $ws.Range("E32") "  +4.49%  "

Set-TextValue $ws.Range("D33") "0.7656"
Set-TextValue $ws.Range("E33") "  +4.46%  "

Set-TextValue $ws.Range("E34") "  +3.67%  "

Set-TextValue $ws.Range("D35") "4.516"
Set-TextValue $ws.Range("E35") "  +1.84%  "

Set-TextValue $ws.Range("D36") "2.709"
Set-TextValue $ws.Range("E36") "  +10.09%  "

Set-TextValue $ws.Range("D37") "1.080"
Set-TextValue $ws.Range("E37") "  +0.86%  "

Set-TextValue $ws.Range("D38") "0.01957"
Set-TextValue $ws.Range("E38") "  +1.03%  "

Set-TextValue $ws.Range("B39") "MXToken"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D39") "2.990"
Set-TextValue $ws.Range("E39") "  +1.63%  "

Set-TextValue $ws.Range("B40") "Hedera"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.05243"
Set-TextValue $ws.Range("E40") "  +0.10%  "

Set-TextValue $ws.Range("D41") "7.075"
Set-TextValue $ws.Range("E41") "  -0.81%  "

Set-TextValue $ws.Range("D42") "0.5146"
Set-TextValue $ws.Range("E42") "  -0.07%  "

Set-TextValue $ws.Range("E43") "  +1.03%  "

Set-TextValue $ws.Range("D44") "8.415"
Set-TextValue $ws.Range("E44") "  +2.39%  "

Set-TextValue $ws.Range("D45") "0.4811"
Set-TextValue $ws.Range("E45") "  -0.28%  "

Set-TextValue $ws.Range("D46") "10.37"
Set-TextValue $ws.Range("E46") "  +1.90%  "

Set-TextValue $ws.Range("D47") "1.002"

Set-TextValue $ws.Range("D48") "102.79"
Set-TextValue $ws.Range("E48") "  +0.64%  "

Set-TextValue $ws.Range("D49") "1.647"
Set-TextValue $ws.Range("E49") "  +1.40%  "

Set-TextValue $ws.Range("D50") "0.06221"
Set-TextValue $ws.Range("E50") "  -0.02%  "

Set-TextValue $ws.Range("D51") "65.44"
Set-TextValue $ws.Range("E51") "  +1.66%  "
